$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 826.7586
$ws.Range("I53").Value = 789.17645
$ws.Range("J53").Value = 880
$ws.Range("K53").Value = 789.17645
$ws.Range("L53").Value = 880
$ws.Range("M53").Value = -152.17645
$ws.Range("N53").Value = -2154

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 238.5
$ws.Range("I92").Value = 238.5
$ws.Range("K92").Value = 238.5
$ws.Range("M92").Value = 1009.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2551.7368
$ws.Range("I141").Value = 2667.2354
$ws.Range("K141").Value = 8001.706200000001
$ws.Range("M141").Value = -2821.706200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4166.875
$ws.Range("I61").Value = 3381.5715
$ws.Range("J61").Value = 4777.6665
$ws.Range("K61").Value = 3381.5715
$ws.Range("L61").Value = 4777.6665
$ws.Range("M61").Value = -3169.5715
$ws.Range("N61").Value = -5201.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4166.875
$ws.Range("I136").Value = 3381.5715
$ws.Range("J136").Value = 4777.6665
$ws.Range("K136").Value = 10144.7145
$ws.Range("L136").Value = 14332.9995
$ws.Range("M136").Value = -7594.7145
$ws.Range("N136").Value = -19432.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1992.6666
$ws.Range("I99").Value = 1992.6666
$ws.Range("K99").Value = 1992.6666
$ws.Range("M99").Value = -494.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3260.0715
$ws.Range("I107").Value = 1290.5
$ws.Range("K107").Value = 1290.5
$ws.Range("M107").Value = 629.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 85725.664
$ws.Range("J132").Value = 85725.664
$ws.Range("L132").Value = 85725.664
$ws.Range("N132").Value = -95845.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 99260
$ws.Range("J140").Value = 99260
$ws.Range("L140").Value = 99260
$ws.Range("N140").Value = -109620

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1794.8
$ws.Range("J31").Value = 2098
$ws.Range("L31").Value = 2098
$ws.Range("N31").Value = -2688

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1794.8
$ws.Range("J34").Value = 2098
$ws.Range("L34").Value = 2098
$ws.Range("N34").Value = -2502

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 69000
$ws.Range("J98").Value = 69000
$ws.Range("L98").Value = 69000
$ws.Range("N98").Value = -73492

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5379.174
$ws.Range("I132").Value = 5172.1
$ws.Range("K132").Value = 15516.3
$ws.Range("M132").Value = -12986.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 75
$ws.Range("J2").Value = 75
$ws.Range("L2").Value = 450
$ws.Range("N2").Value = -676

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 234.3125
$ws.Range("J33").Value = 234.3125
$ws.Range("L33").Value = 1405.875
$ws.Range("N33").Value = -1971.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1747.5
$ws.Range("I51").Value = 2000
$ws.Range("K51").Value = 6000
$ws.Range("M51").Value = -5540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 1767.1428
$ws.Range("I61").Value = 3913.3333
$ws.Range("J61").Value = 157.5
$ws.Range("K61").Value = 11739.9999
$ws.Range("L61").Value = 472.5
$ws.Range("M61").Value = -11524.9999
$ws.Range("N61").Value = -902.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1441.7142
$ws.Range("I68").Value = 1348.6666
$ws.Range("K68").Value = 4045.9998
$ws.Range("M68").Value = -3234.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 3002.2
$ws.Range("J69").Value = 3249.75
$ws.Range("L69").Value = 9749.25
$ws.Range("N69").Value = -11371.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1441.7142
$ws.Range("I71").Value = 1348.6666
$ws.Range("K71").Value = 12137.9994
$ws.Range("M71").Value = -8081.999400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 3002.2
$ws.Range("J72").Value = 3249.75
$ws.Range("L72").Value = 29247.75
$ws.Range("N72").Value = -37359.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1368.6072
$ws.Range("J131").Value = 1731
$ws.Range("L131").Value = 5193
$ws.Range("N131").Value = -15273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 12024
$ws.Range("I57").Value = 5950
$ws.Range("K57").Value = 5950
$ws.Range("M57").Value = -5130

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10490.6
$ws.Range("J70").Value = 10658.667
$ws.Range("L70").Value = 10658.667
$ws.Range("N70").Value = -11198.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 10490.6
$ws.Range("J73").Value = 10658.667
$ws.Range("L73").Value = 10658.667
$ws.Range("N73").Value = -12530.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7821.222
$ws.Range("I80").Value = 3398.6
$ws.Range("J80").Value = 13349.5
$ws.Range("K80").Value = 3398.6
$ws.Range("L80").Value = 13349.5
$ws.Range("M80").Value = -2400.6
$ws.Range("N80").Value = -15345.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 7821.222
$ws.Range("I83").Value = 3398.6
$ws.Range("J83").Value = 13349.5
$ws.Range("K83").Value = 16993
$ws.Range("L83").Value = 66747.5
$ws.Range("M83").Value = -12001
$ws.Range("N83").Value = -76731.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2765.5
$ws.Range("I122").Value = 2874.7
$ws.Range("K122").Value = 8624.099999999999
$ws.Range("M122").Value = -6174.099999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2759.25
$ws.Range("I7").Value = 2759.25
$ws.Range("K7").Value = 2759.25
$ws.Range("M7").Value = -2647.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2759.25
$ws.Range("I126").Value = 2759.25
$ws.Range("K126").Value = 8277.75
$ws.Range("M126").Value = -5807.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3299.875
$ws.Range("I62").Value = 2699.5
$ws.Range("K62").Value = 2699.5
$ws.Range("M62").Value = -2075.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 3299.875
$ws.Range("I65").Value = 2699.5
$ws.Range("K65").Value = 13497.5
$ws.Range("M65").Value = -10377.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2814.8235
$ws.Range("I96").Value = 3259.2727
$ws.Range("K96").Value = 3259.2727
$ws.Range("M96").Value = -1886.2727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1142.5714
$ws.Range("I132").Value = 1083
$ws.Range("K132").Value = 3249
$ws.Range("M132").Value = -719
